$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2895.4614
$ws.Range("I33").Value = 2940.1667
$ws.Range("J33").Value = 2857.1428
$ws.Range("K33").Value = 2940.1667
$ws.Range("L33").Value = 2857.1428
$ws.Range("M33").Value = -2711.1667
$ws.Range("N33").Value = -3315.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 231.05263
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1496.1177
$ws.Range("I43").Value = 870
$ws.Range("J43").Value = 1837.6364
$ws.Range("K43").Value = 870
$ws.Range("L43").Value = 1837.6364
$ws.Range("M43").Value = -801
$ws.Range("N43").Value = -1975.6364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1702
$ws.Range("I98").Value = 1539.75
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 1539.75
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -41.75
$ws.Range("N98").Value = -5996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1702
$ws.Range("I122").Value = 1539.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4619.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2169.25
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4524.1333
$ws.Range("I137").Value = 4704.4287
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 14113.2861
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -11563.2861
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 190038.77
$ws.Range("I138").Value = 4251.778
$ws.Range("J138").Value = 224162.89
$ws.Range("K138").Value = 12755.334
$ws.Range("L138").Value = 672488.67
$ws.Range("M138").Value = -7615.334000000001
$ws.Range("N138").Value = -682768.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3410.2104
$ws.Range("I45").Value = 2986.75
$ws.Range("K45").Value = 2986.75
$ws.Range("M45").Value = -2609.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1374.7084
$ws.Range("I74").Value = 1045.5625
$ws.Range("J74").Value = 2033
$ws.Range("K74").Value = 1045.5625
$ws.Range("L74").Value = 2033
$ws.Range("M74").Value = -171.5625
$ws.Range("N74").Value = -3781

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1374.7084
$ws.Range("I77").Value = 1045.5625
$ws.Range("J77").Value = 2033
$ws.Range("K77").Value = 5227.8125
$ws.Range("L77").Value = 10165
$ws.Range("M77").Value = -859.8125
$ws.Range("N77").Value = -18901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 44559.223
$ws.Range("J112").Value = 44559.223
$ws.Range("L112").Value = 44559.223
$ws.Range("N112").Value = -47513.223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 33328
$ws.Range("J123").Value = 33328
$ws.Range("L123").Value = 33328
$ws.Range("N123").Value = -43128

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4443.25
$ws.Range("I132").Value = 3751.4211
$ws.Range("J132").Value = 5454.385
$ws.Range("K132").Value = 11254.2633
$ws.Range("L132").Value = 16363.155
$ws.Range("M132").Value = -8724.263300000001
$ws.Range("N132").Value = -21423.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 520.96
$ws.Range("I94").Value = 461.75
$ws.Range("J94").Value = 757.8
$ws.Range("K94").Value = 461.75
$ws.Range("L94").Value = 757.8
$ws.Range("M94").Value = -10.75
$ws.Range("N94").Value = -1659.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 35000
$ws.Range("J112").Value = 35000
$ws.Range("L112").Value = 35000
$ws.Range("N112").Value = -37954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4883.091
$ws.Range("I134").Value = 7633.3335
$ws.Range("J134").Value = 3851.75
$ws.Range("K134").Value = 22900.0005
$ws.Range("L134").Value = 11555.25
$ws.Range("M134").Value = -20365.0005
$ws.Range("N134").Value = -16625.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 499.5
$ws.Range("I16").Value = 499.66666
$ws.Range("J16").Value = 499
$ws.Range("K16").Value = 499.66666
$ws.Range("L16").Value = 499
$ws.Range("M16").Value = -212.66666
$ws.Range("N16").Value = -1073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2791.4062
$ws.Range("I31").Value = 1021.13043
$ws.Range("K31").Value = 1021.13043
$ws.Range("M31").Value = -726.13043

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2791.4062
$ws.Range("I34").Value = 1021.13043
$ws.Range("K34").Value = 1021.13043
$ws.Range("M34").Value = -819.13043

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 499.5
$ws.Range("I113").Value = 499.66666
$ws.Range("J113").Value = 499
$ws.Range("K113").Value = 499.66666
$ws.Range("L113").Value = 499
$ws.Range("M113").Value = 1670.33334
$ws.Range("N113").Value = -4839

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 20833592
$ws.Range("I50").Value = 245.45454
$ws.Range("J50").Value = 66666950
$ws.Range("K50").Value = 736.3636200000001
$ws.Range("L50").Value = 200000850
$ws.Range("M50").Value = -255.3636200000001
$ws.Range("N50").Value = -200001812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 20833592
$ws.Range("I53").Value = 245.45454
$ws.Range("J53").Value = 66666950
$ws.Range("K53").Value = 736.3636200000001
$ws.Range("L53").Value = 200000850
$ws.Range("M53").Value = -255.3636200000001
$ws.Range("N53").Value = -200001812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 944.34375
$ws.Range("J131").Value = 1038.8928
$ws.Range("L131").Value = 3116.6784
$ws.Range("N131").Value = -13196.6784

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3577.5312
$ws.Range("I132").Value = 2851.6191
$ws.Range("J132").Value = 3932.0466
$ws.Range("K132").Value = 25664.5719
$ws.Range("L132").Value = 35388.4194
$ws.Range("M132").Value = -23134.5719
$ws.Range("N132").Value = -40448.4194

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1616.1364
$ws.Range("I140").Value = 899.0625
$ws.Range("K140").Value = 2697.1875
$ws.Range("M140").Value = 2482.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2833.3333
$ws.Range("I102").Value = 2750
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2750
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1128
$ws.Range("N102").Value = -6244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 30833.334
$ws.Range("J111").Value = 30833.334
$ws.Range("L111").Value = 30833.334
$ws.Range("N111").Value = -36967.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1968.6666
$ws.Range("J113").Value = 1968.6666
$ws.Range("L113").Value = 1968.6666
$ws.Range("N113").Value = -6308.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 111115660
$ws.Range("I7").Value = 125004250
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 125004250
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -125004138
$ws.Range("N7").Value = -7224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10801.182
$ws.Range("I93").Value = 21100.2
$ws.Range("J93").Value = 2218.6667
$ws.Range("K93").Value = 21100.2
$ws.Range("L93").Value = 2218.6667
$ws.Range("M93").Value = -19852.2
$ws.Range("N93").Value = -4714.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 20600
$ws.Range("J103").Value = 20600
$ws.Range("L103").Value = 20600
$ws.Range("N103").Value = -22944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 29724.875
$ws.Range("J110").Value = 29724.875
$ws.Range("L110").Value = 29724.875
$ws.Range("N110").Value = -37904.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 111115660
$ws.Range("I126").Value = 125004250
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 375012750
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -375010280
$ws.Range("N126").Value = -25940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1252.2778
$ws.Range("I126").Value = 1088.7273
$ws.Range("J126").Value = 1509.2858
$ws.Range("K126").Value = 3266.1819
$ws.Range("L126").Value = 4527.857400000001
$ws.Range("M126").Value = -796.1819
$ws.Range("N126").Value = -9467.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3970992.8
$ws.Range("I132").Value = 3781.875
$ws.Range("J132").Value = 6412353.5
$ws.Range("K132").Value = 11345.625
$ws.Range("L132").Value = 19237060.5
$ws.Range("M132").Value = -8815.625
$ws.Range("N132").Value = -19242120.5
